# Update the "9-supply" sheet's num_lines table (existing + planned TDP lines)
# and touch the data range's formatting (mirrors selecting the block and
# clearing/re-applying fill), matching the author's commit:
# "updated num_lines based on TDP existing and planned lines"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("9-supply")

# --- Data updates (num_parallel_275 / num_parallel_400 / num_parallel_765) ---

# row 2: PINETOWN(idx0)/WESTERN CAPE-NORTHERN CAPE corridor
$ws.Range("F2").Value = 4

# row 4
$ws.Range("A4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2

# row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 2

# row 7
$ws.Range("A7").Value = 2
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0

# row 8
$ws.Range("D8").Value = 2

# row 9
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 1

# row 10
$ws.Range("F10").Value = 0

# row 13
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 3

# row 14
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 1

# row 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1

# row 16
$ws.Range("E16").Value = 1

# row 18
$ws.Range("E18").Value = 1

# --- Formatting touch over the whole data body (A2:F18) ---
# (selecting the block and clicking "No Fill" leaves an explicit,
# visually-unchanged style on every cell in the range)
$body = $ws.Range("A2:F18")
$body.Interior.ColorIndex = -4142

# --- Selection left on row 17 (whole row), as in the saved file ---
$ws.Range("A17:XFD17").Select()
